# Generalização do RDM e escolha da estratégia finalizada!
#
# - Widen the "Taxa de Contatos" (aContactRate) Min/Max bounds on the
#   "params" sheet (row 3: C3 Min 1 -> 20, D3 Max 50 -> 100).
# - Move the active selection on "levers" from C18 to the D2:D18 range.
# - Switch the active/selected worksheet from "levers" to "params", and
#   move its selection to A4.

$wb = $excel.ActiveWorkbook

$wsParams = $wb.Worksheets.Item("params")
$wsLevers = $wb.Worksheets.Item("levers")

# Update the Min/Max bounds for the Contact Rate lever on "params".
$wsParams.Range("C3").Value = 20
$wsParams.Range("D3").Value = 100

# Update the selection on "levers" to D2:D18 (no longer the active tab).
$wsLevers.Range("D2:D18").Select()

# Make "params" the active sheet and select A4 there.
$wsParams.Activate()
$wsParams.Range("A4").Select()
